# Insert a new daily price record as row 4 on the single data sheet.
# This pushes the former rows 4..79 down to 5..80 (dimension grows from
# A1:T79 to A1:T80), and the new row 4 gets a fresh set of observations
# (columns A-L, T repeat the constant/categorical values used throughout
# the sheet; D, M, N, O, P, Q, R, S hold the new record's data).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing rows 4-79 down to 5-80, creating a blank row 4.
$ws.Rows.Item(4).Insert()

# Populate the new row 4 with the inserted record.
$ws.Range("A4").Value = 5
$ws.Range("B4").Value = "Macroferia Regional de Talca"
$ws.Range("C4").Value = "Maule"
$ws.Range("D4").Value = 44882
$ws.Range("E4").Value = 7
$ws.Range("F4").Value = "Fruta"
$ws.Range("G4").Value = 100101
$ws.Range("H4").Value = "Berries"
$ws.Range("I4").Value = 100101001
$ws.Range("J4").Value = "Arándano (blue)"
$ws.Range("K4").Value = "Sin especificar"
$ws.Range("L4").Value = "Primera"
$ws.Range("M4").Value = 170
$ws.Range("N4").Value = 6000
$ws.Range("O4").Value = 6000
$ws.Range("P4").Value = 6000
$ws.Range("Q4").Value = "`$/bandeja 2 kilos"
$ws.Range("R4").Value = "Región de O'Higgins"
$ws.Range("S4").Value = 3000
$ws.Range("T4").Value = 2
